# "cambio columnas excel google"
# Rename the "link" and "status" header columns on the Google product feed
# template sheet:
#   D1: "enlace al producto" -> "link"
#   E1: "estado (nuevo)"     -> "estado"
# and leave the cursor on E2, matching where the editor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "link"
$ws.Range("E1").Value = "estado"

$ws.Range("E2").Select()
